$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '42.718.85'
$ws.Range("E2").Value = '  -0.35%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.545.36'
$ws.Range("E3").Value = '  +0.26%  '
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '309.90'
$ws.Range("E5").Value = '  -2.70%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '99.81'
$ws.Range("E6").Value = '  +2.72%  '
$ws.Range("E7").Value = '  -0.83%  '
$ws.Range("E8").Value = '  -0.04%  '
$ws.Range("E9").Value = '  -0.45%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '35.90'
$ws.Range("E10").Value = '  -0.77%  '
$ws.Range("E11").Value = '  -1.65%  '
$ws.Range("E12").Value = '  -2.79%  '
$ws.Range("E13").Value = '  -0.94%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '2.937.63'
$ws.Range("E14").Value = '  +0.18%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '15.99'
$ws.Range("E15").Value = '  +5.31%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '2.549.84'
$ws.Range("E16").Value = '  -2.52%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.840'
$ws.Range("E17").Value = '  -1.60%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '42.716.96'
$ws.Range("E18").Value = '  -0.58%  '
$ws.Range("E19").Value = '  -2.04%  '
$ws.Range("B20").Value = 'ShibaInu'
$ws.Range("C20").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.0₃0955'
$ws.Range("E20").Value = '  -1.56%  '
$ws.Range("B21").Value = 'InternetComputer(DFINITY)'
$ws.Range("C21").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '12.36'
$ws.Range("E21").Value = '  -3.46%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '69.42'
$ws.Range("E22").Value = '  -0.57%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '247.91'
$ws.Range("E23").Value = '  -2.33%  '
$ws.Range("E24").Value = '  -2.23%  '
$ws.Range("E25").Value = '  -0.02%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '26.57'
$ws.Range("E26").Value = '  +0.20%  '
$ws.Range("E27").Value = '  -0.02%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.37'
$ws.Range("E28").Value = '  -1.76%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '40.01'
$ws.Range("E29").Value = '  -2.34%  '
$ws.Range("E30").Value = '  -3.38%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '158.60'
$ws.Range("E31").Value = '  +0.41%  '
$ws.Range("E32").Value = '  -2.99%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0806'
$ws.Range("E33").Value = '  +1.42%  '
$ws.Range("E34").Value = '  -1.92%  '
$ws.Range("E35").Value = '  -3.81%  '
$ws.Range("E36").Value = '  -3.49%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.62'
$ws.Range("E37").Value = '  +5.62%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '18.31'
$ws.Range("E38").Value = '  -4.64%  '
$ws.Range("E39").Value = '  -1.51%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.117'
$ws.Range("E40").Value = '  -0.94%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '22.36'
$ws.Range("E41").Value = '  +0.80%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '4.16'
$ws.Range("E42").Value = '  +8.48%  '
$ws.Range("E43").Value = '  -0.15%  '
$ws.Range("E44").Value = '  -1.52%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '3.24'
$ws.Range("E45").Value = '  -1.35%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.991.08'
$ws.Range("E46").Value = '  -1.24%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '9.04'
$ws.Range("E47").Value = '  -1.66%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.784.21'
$ws.Range("E48").Value = '  -0.09%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '81.25'
$ws.Range("E49").Value = '  -3.79%  '
$ws.Range("E50").Value = '  +0.23%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '73.37'
$ws.Range("E51").Value = '  -3.67%  '
